$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old EditorUI row ("editor menu" / "[Space] - Show menu") at row 83.
$ws.Rows(83).Delete()
